# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates (and a few cell additions/removals) to the
# "Leve Profit" sheets, per the target diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 268
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("H39").Value = 608
$ws.Range("I39").Value = 78.27778000000001
$ws.Range("J39").Value = 3786.3333
$ws.Range("K39").Value = 234.83334
$ws.Range("L39").Value = 11358.9999
$ws.Range("M39").Value = 61.16665999999998
$ws.Range("N39").Value = -11950.9999
$ws.Range("H51").Value = 36667.332
$ws.Range("J51").Value = 36667.332
$ws.Range("L51").Value = 36667.332
$ws.Range("N51").Value = -37635.332
$ws.Range("H111").Value = 528.375
$ws.Range("I111").Value = 297
$ws.Range("J111").Value = 759.75
$ws.Range("K111").Value = 891
$ws.Range("L111").Value = 2279.25
$ws.Range("M111").Value = 2176
$ws.Range("N111").Value = -8413.25
$ws.Range("H139").Value = 166988.6
$ws.Range("J139").Value = 166988.6
$ws.Range("L139").Value = 166988.6
$ws.Range("N139").Value = -177268.6
$ws.Range("N2").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 979576.9399999999
$ws.Range("I2").Value = 1305202.8
$ws.Range("J2").Value = 2699.5
$ws.Range("K2").Value = 1305202.8
$ws.Range("L2").Value = 2699.5
$ws.Range("M2").Value = -1305089.8
$ws.Range("N2").Value = -2925.5
$ws.Range("H32").Value = 8663.143
$ws.Range("I32").Value = 8663.143
$ws.Range("K32").Value = 8663.143
$ws.Range("M32").Value = -8376.143
$ws.Range("H45").Value = 3187.2856
$ws.Range("I45").Value = 3187.2856
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3187.2856
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2810.2856
$ws.Range("H97").Value = 2034.5714
$ws.Range("I97").Value = 2186.25
$ws.Range("J97").Value = 1832.3334
$ws.Range("K97").Value = 2186.25
$ws.Range("L97").Value = 1832.3334
$ws.Range("M97").Value = -1690.25
$ws.Range("N97").Value = -2824.3334
$ws.Range("H116").Value = 979576.9399999999
$ws.Range("I116").Value = 1305202.8
$ws.Range("J116").Value = 2699.5
$ws.Range("K116").Value = 1305202.8
$ws.Range("L116").Value = 2699.5
$ws.Range("M116").Value = -1302908.8
$ws.Range("N116").Value = -7287.5
$ws.Range("N45").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 979576.9399999999
$ws.Range("I3").Value = 1305202.8
$ws.Range("J3").Value = 2699.5
$ws.Range("K3").Value = 1305202.8
$ws.Range("L3").Value = 2699.5
$ws.Range("M3").Value = -1305088.8
$ws.Range("N3").Value = -2927.5
$ws.Range("H94").Value = 769.5909
$ws.Range("I94").Value = 787.1905
$ws.Range("J94").Value = 400
$ws.Range("K94").Value = 787.1905
$ws.Range("L94").Value = 400
$ws.Range("M94").Value = -336.1905
$ws.Range("N94").Value = -1302
$ws.Range("H99").Value = 2023.4615
$ws.Range("I99").Value = 1781
$ws.Range("K99").Value = 1781
$ws.Range("M99").Value = -283
$ws.Range("H107").Value = 1042.5714
$ws.Range("I107").Value = 984.57574
$ws.Range("K107").Value = 984.57574
$ws.Range("M107").Value = 935.42426
$ws.Range("H134").Value = 4697.727
$ws.Range("I134").Value = 4567
$ws.Range("K134").Value = 13701
$ws.Range("M134").Value = -11166

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9238.588
$ws.Range("I58").Value = 10303.786
$ws.Range("J58").Value = 4267.6665
$ws.Range("K58").Value = 10303.786
$ws.Range("L58").Value = 4267.6665
$ws.Range("M58").Value = -10100.786
$ws.Range("N58").Value = -4673.6665
$ws.Range("H99").Value = 5844
$ws.Range("I99").Value = 5441
$ws.Range("J99").Value = 6247
$ws.Range("K99").Value = 5441
$ws.Range("L99").Value = 6247
$ws.Range("M99").Value = -3943
$ws.Range("N99").Value = -9243
$ws.Range("H107").Value = 866106.7
$ws.Range("I107").Value = 1276285.2
$ws.Range("K107").Value = 1276285.2
$ws.Range("M107").Value = -1274365.2
$ws.Range("H126").Value = 5844
$ws.Range("I126").Value = 5441
$ws.Range("J126").Value = 6247
$ws.Range("K126").Value = 16323
$ws.Range("L126").Value = 18741
$ws.Range("M126").Value = -13853
$ws.Range("N126").Value = -23681
$ws.Range("H130").Value = 24374.25
$ws.Range("J130").Value = 24374.25
$ws.Range("L130").Value = 24374.25
$ws.Range("N130").Value = -34414.25
$ws.Range("H136").Value = 9238.588
$ws.Range("I136").Value = 10303.786
$ws.Range("J136").Value = 4267.6665
$ws.Range("K136").Value = 30911.358
$ws.Range("L136").Value = 12802.9995
$ws.Range("M136").Value = -28361.358
$ws.Range("N136").Value = -17902.9995
$ws.Range("H141").Value = 591654.5600000001
$ws.Range("J141").Value = 633350
$ws.Range("L141").Value = 633350
$ws.Range("N141").Value = -643710

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 11387352
$ws.Range("I7").Value = 13284578
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 39853734
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = -39853622
$ws.Range("N7").Value = -12224
$ws.Range("H8").Value = 349.41666
$ws.Range("I8").Value = 349.41666
$ws.Range("K8").Value = 1048.24998
$ws.Range("M8").Value = -909.2499800000001
$ws.Range("H36").Value = 1483
$ws.Range("I36").Value = 1483
$ws.Range("K36").Value = 4449
$ws.Range("M36").Value = -4280
$ws.Range("H38").Value = 325.125
$ws.Range("I38").Value = 316.66666
$ws.Range("J38").Value = 333.58334
$ws.Range("K38").Value = 949.9999799999999
$ws.Range("L38").Value = 1000.75002
$ws.Range("M38").Value = -602.9999799999999
$ws.Range("N38").Value = -1694.75002
$ws.Range("H40").Value = 81.375
$ws.Range("I40").Value = 81.375
$ws.Range("K40").Value = 325.5
$ws.Range("M40").Value = -256.5
$ws.Range("H68").Value = 1659
$ws.Range("J68").Value = 1754.875
$ws.Range("L68").Value = 5264.625
$ws.Range("N68").Value = -6886.625
$ws.Range("H71").Value = 1659
$ws.Range("J71").Value = 1754.875
$ws.Range("L71").Value = 15793.875
$ws.Range("N71").Value = -23905.875
$ws.Range("H92").Value = 458.58334
$ws.Range("I92").Value = 458.58334
$ws.Range("K92").Value = 1375.75002
$ws.Range("M92").Value = -127.7500199999999
$ws.Range("H108").Value = 2356.5715
$ws.Range("I108").Value = 1716
$ws.Range("K108").Value = 5148
$ws.Range("M108").Value = -2268
$ws.Range("H113").Value = 68395.336
$ws.Range("J113").Value = 1839.1428
$ws.Range("L113").Value = 5517.428400000001
$ws.Range("N113").Value = -9857.428400000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 26892.666
$ws.Range("J15").Value = 27629.25
$ws.Range("L15").Value = 27629.25
$ws.Range("N15").Value = -28205.25
$ws.Range("H33").Value = 29945
$ws.Range("I33").Value = 29945
$ws.Range("K33").Value = 29945
$ws.Range("M33").Value = -29693
$ws.Range("H81").Value = 26892.666
$ws.Range("J81").Value = 27629.25
$ws.Range("L81").Value = 27629.25
$ws.Range("N81").Value = -29625.25
$ws.Range("H84").Value = 26892.666
$ws.Range("J84").Value = 27629.25
$ws.Range("L84").Value = 82887.75
$ws.Range("N84").Value = -92871.75
$ws.Range("H97").Value = 1248.1666
$ws.Range("I97").Value = 997.25
$ws.Range("K97").Value = 997.25
$ws.Range("M97").Value = -501.25
$ws.Range("H105").Value = 37932.332
$ws.Range("J105").Value = 69999
$ws.Range("L105").Value = 69999
$ws.Range("N105").Value = -76987
$ws.Range("H132").Value = 8659.299999999999
$ws.Range("I132").Value = 8600.333000000001
$ws.Range("K132").Value = 25800.999
$ws.Range("M132").Value = -23270.999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1749.9667
$ws.Range("I22").Value = 1125.1364
$ws.Range("J22").Value = 3468.25
$ws.Range("K22").Value = 1125.1364
$ws.Range("L22").Value = 3468.25
$ws.Range("M22").Value = -830.1364000000001
$ws.Range("N22").Value = -4058.25
$ws.Range("H27").Value = 1749.9667
$ws.Range("I27").Value = 1125.1364
$ws.Range("J27").Value = 3468.25
$ws.Range("K27").Value = 1125.1364
$ws.Range("L27").Value = 3468.25
$ws.Range("M27").Value = -1018.1364
$ws.Range("N27").Value = -3682.25
$ws.Range("H46").Value = 1789.75
$ws.Range("I46").Value = 1609.8889
$ws.Range("J46").Value = 2329.3333
$ws.Range("K46").Value = 1609.8889
$ws.Range("L46").Value = 2329.3333
$ws.Range("M46").Value = -1421.8889
$ws.Range("N46").Value = -2705.3333
$ws.Range("H93").Value = 2457.2666
$ws.Range("I93").Value = 2241.3635
$ws.Range("J93").Value = 3051
$ws.Range("K93").Value = 2241.3635
$ws.Range("L93").Value = 3051
$ws.Range("M93").Value = -993.3634999999999
$ws.Range("N93").Value = -5547

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 84949.5
$ws.Range("J42").Value = 84949.5
$ws.Range("L42").Value = 84949.5
$ws.Range("N42").Value = -85705.5
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("H132").Value = 4372.517
$ws.Range("I132").Value = 5233.8823
$ws.Range("K132").Value = 15701.6469
$ws.Range("M132").Value = -13171.6469
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H136").Value = 2617.4375
$ws.Range("I136").Value = 2639.6924
$ws.Range("K136").Value = 7919.0772
$ws.Range("M136").Value = -5369.0772
$ws.Range("N46").ClearContents()
$ws.Range("N134").ClearContents()

Write-Output "Applied 242 cell updates and 4 cell clears across 8 sheets."
